$d = $word.ActiveDocument

# 1) Split the title run "ASTERA System " into "MUCKE" + " System "
#    (same bold/lang formatting on both resulting runs).
$titleRun = $d.Range(0, 6)
$titleRun.Font.Bold = $false
$titleRun.Text = "MUCKE"
$afterMucke = $d.Range(0, 5)
$afterMucke.Font.Bold = $true

# 2) Move the "_GoBack" bookmark from the empty paragraph right before the
#    "Configuration Mechanism" heading to the empty paragraph right after
#    the "One high-level paragraph..." text (Adding a bookmark with an
#    existing name relocates it).
$target = $d.Paragraphs(4).Range
$d.Bookmarks.Add("_GoBack", $target)
